$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the literal "<br/>" markers inside these cell values with real
# line breaks (newline characters), matching the diff against
# sharedStrings.xml. Cells whose text begins with "=" are written with a
# leading apostrophe (Excel's "treat as text" prefix) so the engine does
# not try to parse the value as a formula; the apostrophe itself is not
# part of the stored text. The cell style is then reset to "Normal" so the
# quote-prefix flag introduced on the style doesn't linger, keeping the
# resulting formatting unchanged relative to the original.

function Set-PlainText($cell, $text) {
    $range = $ws.Range($cell)
    if ($text.StartsWith("=")) {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-PlainText "B24" "=<19 Manufacturing`n=<49 Service"
Set-PlainText "D24" "=< MNT 250 Millionlion Manufacturing`n=< MNT 1 Billionlion Service"

Set-PlainText "B25" "=<149 Wholesale trade`n=<199 Retail Trade`n=<199 Manufacturing"
Set-PlainText "D25" "=< MNT 1.5 Billionlion Wholesale trade`n=< MNT 1.5 Billionlion Retail Trade`n=< MNT 1.5 Billionlion Manufacturing"

Set-PlainText "B26" ">149 Wholesale trade`n>199 Retail Trade`n>199 Manufacturing"
Set-PlainText "D26" "> MNT 1.5 Billionlion Wholesale trade`n> MNT 1.5 Billionlion Retail Trade`n> MNT 1.5 Billionlion Manufacturing"
